$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.144.20"
$ws.Range("E2").Value = "  -0.28%  "

$ws.Range("D3").Value = "1.841.72"
$ws.Range("E3").Value = "  -0.42%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6868"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9994"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3016"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.82%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07451"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.53%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.30%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07649"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.29%  "

$ws.Range("D12").Value = "1.838.99"
$ws.Range("E12").Value = "  -0.44%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.060"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6828"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.70%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "87.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.82%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.172"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.71%  "

$ws.Range("D17").Value = "29.131.26"
$ws.Range("E17").Value = "  -0.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008157"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.03%  "

$ws.Range("D19").Value = "2.078.08"
$ws.Range("E19").Value = "  -0.79%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "228.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.82%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9996"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.393"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.81%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9995"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1452"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.32%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.762"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.95%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.510"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.274"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.97%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.140"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.194"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.85%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05247"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.43%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7589"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.851"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.85%  "

$ws.Range("E36").Value = "  -1.31%  "

$ws.Range("E37").Value = "  -0.43%  "

$ws.Range("D38").Value = "1.302.97"
$ws.Range("E38").Value = "  -1.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01835"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.724"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.26%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9295"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.69%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.916"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "104.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9991"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.08%  "

$ws.Range("D45").Value = "1.980.34"
$ws.Range("E45").Value = "  -0.57%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5197"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.30%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000123"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.27%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.531"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.84%  "

$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.770"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.12%  "

$ws.Range("B51").Value = "XinFinNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07415"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +17.50%  "
